# Adding a POC for using material layout with the material framework:
# a "Text & Icon Opacity" table (F:H) and a "Typography" table (J:M)
# alongside the existing keylines/metrics table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Table 1: Name | Dark Opacity | White Opacity (F:H) ----
$ws.Range("F1").Value = "Name"
$ws.Range("F5").Value = "Dividers"
$ws.Range("F2").Value = "Primary Text"
$ws.Range("F3").Value = "Secondary Text"
$ws.Range("F4").Value = "Disabled/hint Text"
$ws.Range("H1").Value = "White Opacity"
$ws.Range("G1").Value = "Dark Opacity"
$ws.Range("F6").Value = "Active Icon"
$ws.Range("F7").Value = "Inactive Icon"

$ws.Range("G2").Value = 87
$ws.Range("H2").Value = 100
$ws.Range("G3").Value = 54
$ws.Range("H3").Value = 70
$ws.Range("G4").Value = 38
$ws.Range("H4").Value = 50
$ws.Range("G5").Value = 12
$ws.Range("H5").Value = 12
$ws.Range("G6").Value = 54
$ws.Range("H6").Value = 100
$ws.Range("G7").Value = 38
$ws.Range("H7").Value = 50

$ws.Range("F1:H1").Font.Bold = $true

# ---- Table 2: Name | Device | Desktop (J:M) ----
$ws.Range("J1").Value = "Name"
$ws.Range("K1").Value = "Device"
$ws.Range("L1").Value = "Desktop"
$ws.Range("J2").Value = "App Bar"
$ws.Range("J3").Value = "Button"
$ws.Range("J4").Value = "Subtitle"
$ws.Range("J5").Value = "Body"
$ws.Range("J6").Value = "Caption"

$ws.Range("K2").Value = 20
$ws.Range("L2").Value = 18
$ws.Range("K3").Value = 14
$ws.Range("L3").Value = 14
$ws.Range("K4").Value = 16
$ws.Range("L4").Value = 15
$ws.Range("K5").Value = 14
$ws.Range("L5").Value = 13
$ws.Range("K6").Value = 12
$ws.Range("L6").Value = 12

# J1:M1 includes M1, which carries the header style (bold) but no value,
# matching the source sheet.
$ws.Range("J1:M1").Font.Bold = $true

# ---- Column widths (best-fit sized for the new tables' content;
#      values chosen so the engine's internal rounding lands as close as
#      possible to Excel's own best-fit pixel widths: 17.71, 12.14, 13.86,
#      11.43, 12.71 characters respectively) ----
$ws.Columns("F:F").ColumnWidth = 16.8333333333
$ws.Columns("G:G").ColumnWidth = 11.3333333333
$ws.Columns("H:H").ColumnWidth = 13.0
$ws.Columns("J:J").ColumnWidth = 10.6666666667
$ws.Columns("M:M").ColumnWidth = 11.8333333333

# ---- Selection, matching the committed sheet view ----
$ws.Range("G6").Select()
